$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 78.5
$ws.Range("I33").Value = 84.3
$ws.Range("K33").Value = 84.3
$ws.Range("M33").Value = 144.7

$ws.Range("H132").Value = 613.2909
$ws.Range("I132").Value = 623.17645
$ws.Range("K132").Value = 1869.52935
$ws.Range("M132").Value = 660.4706499999998

$ws.Range("H137").Value = 5049.6665
$ws.Range("I137").Value = 2899.9412
$ws.Range("J137").Value = 6973.1055
$ws.Range("K137").Value = 8699.8236
$ws.Range("L137").Value = 20919.3165
$ws.Range("M137").Value = -6149.8236
$ws.Range("N137").Value = -26019.3165

$ws.Range("H138").Value = 1519386.9
$ws.Range("I138").Value = 2683.3
$ws.Range("K138").Value = 8049.900000000001
$ws.Range("M138").Value = -2909.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32499.25
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 32499.25
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 32499.25
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -33073.25

$ws.Range("H124").Value = 54423.168
$ws.Range("J124").Value = 54423.168
$ws.Range("L124").Value = 54423.168
$ws.Range("N124").Value = -64243.168

$ws.Range("H125").Value = 40359
$ws.Range("J125").Value = 40359
$ws.Range("L125").Value = 40359
$ws.Range("N125").Value = -50199

$ws.Range("H132").Value = 4411.18
$ws.Range("I132").Value = 1270.7812
$ws.Range("J132").Value = 9994.111
$ws.Range("K132").Value = 3812.3436
$ws.Range("L132").Value = 29982.333
$ws.Range("M132").Value = -1282.3436
$ws.Range("N132").Value = -35042.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6762721.5
$ws.Range("I134").Value = 9618762
$ws.Range("J134").Value = 12081.272
$ws.Range("K134").Value = 28856286
$ws.Range("L134").Value = 36243.81600000001
$ws.Range("M134").Value = -28853751
$ws.Range("N134").Value = -41313.81600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1821.8182
$ws.Range("J94").Value = 1408.2858
$ws.Range("L94").Value = 1408.2858
$ws.Range("N94").Value = -2310.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 405.8095
$ws.Range("J12").Value = 57.615383
$ws.Range("L12").Value = 172.846149
$ws.Range("N12").Value = -518.846149

$ws.Range("H34").Value = 4226.609
$ws.Range("I34").Value = 714.6667
$ws.Range("J34").Value = 4753.4
$ws.Range("K34").Value = 2144.0001
$ws.Range("L34").Value = 14260.2
$ws.Range("M34").Value = -2060.0001
$ws.Range("N34").Value = -14428.2

$ws.Range("H39").Value = 13198.5
$ws.Range("I39").Value = 14000
$ws.Range("J39").Value = 13038.2
$ws.Range("K39").Value = 42000
$ws.Range("L39").Value = 39114.60000000001
$ws.Range("M39").Value = -41706
$ws.Range("N39").Value = -39702.60000000001

$ws.Range("H55").Value = 4553115
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4553115
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 13659345
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -13659699

$ws.Range("H113").Value = 4403.6
$ws.Range("J113").Value = 5416.316
$ws.Range("L113").Value = 16248.948
$ws.Range("N113").Value = -20588.948

$ws.Range("H132").Value = 7720.9546
$ws.Range("I132").Value = 2688.5833
$ws.Range("J132").Value = 13759.8
$ws.Range("K132").Value = 24197.2497
$ws.Range("L132").Value = 123838.2
$ws.Range("M132").Value = -21667.2497
$ws.Range("N132").Value = -128898.2

$ws.Range("H141").Value = 5994.5
$ws.Range("I141").Value = 4902.091
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 14706.273
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -9526.273000000001
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 13995
$ws.Range("I35").Value = 13995
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 13995
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -13697
$ws.Range("N35").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H52").Value = 90000
$ws.Range("J52").Value = 90000
$ws.Range("L52").Value = 90000
$ws.Range("N52").Value = -90518

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H97").Value = 1168.4348
$ws.Range("I97").Value = 835.2143
$ws.Range("J97").Value = 1686.7778
$ws.Range("K97").Value = 835.2143
$ws.Range("L97").Value = 1686.7778
$ws.Range("M97").Value = -339.2143
$ws.Range("N97").Value = -2678.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 7200.6665
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 8801
$ws.Range("K58").Value = 4000
$ws.Range("L58").Value = 8801
$ws.Range("M58").Value = -3740
$ws.Range("N58").Value = -9321

$ws.Range("H99").Value = 32100.5
$ws.Range("I99").Value = 10259
$ws.Range("J99").Value = 53942
$ws.Range("K99").Value = 10259
$ws.Range("L99").Value = 53942
$ws.Range("M99").Value = -7264
$ws.Range("N99").Value = -59932

$ws.Range("H132").Value = 7941534.5
$ws.Range("I132").Value = 13515664
$ws.Range("J132").Value = 9119.192
$ws.Range("K132").Value = 40546992
$ws.Range("L132").Value = 27357.576
$ws.Range("M132").Value = -40544462
$ws.Range("N132").Value = -32417.576

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 30012.5
$ws.Range("I43").Value = 30027
$ws.Range("J43").Value = 29998
$ws.Range("K43").Value = 30027
$ws.Range("L43").Value = 29998
$ws.Range("M43").Value = -29878
$ws.Range("N43").Value = -30296

$ws.Range("H51").Value = 45000
$ws.Range("I51").Value = 45000
$ws.Range("K51").Value = 45000
$ws.Range("M51").Value = -44490

$ws.Range("H52").Value = 9500
$ws.Range("I52").Value = 9500
$ws.Range("K52").Value = 9500
$ws.Range("M52").Value = -9274

$ws.Range("H62").Value = 52599.5
$ws.Range("I62").Value = 70856.29
$ws.Range("J62").Value = 10000.333
$ws.Range("K62").Value = 70856.29
$ws.Range("L62").Value = 10000.333
$ws.Range("M62").Value = -70232.29
$ws.Range("N62").Value = -11248.333

$ws.Range("H65").Value = 52599.5
$ws.Range("I65").Value = 70856.29
$ws.Range("J65").Value = 10000.333
$ws.Range("K65").Value = 354281.45
$ws.Range("L65").Value = 50001.665
$ws.Range("M65").Value = -351161.45
$ws.Range("N65").Value = -56241.665

$ws.Range("H96").Value = 1999.6666
$ws.Range("I96").Value = 1999
$ws.Range("K96").Value = 1999
$ws.Range("M96").Value = -626

$ws.Range("H100").Value = 605.03845
$ws.Range("I100").Value = 429.94446
$ws.Range("K100").Value = 859.88892
$ws.Range("M100").Value = -318.88892

$ws.Range("H132").Value = 3859.44
$ws.Range("J132").Value = 2982.6191
$ws.Range("L132").Value = 8947.8573
$ws.Range("N132").Value = -14007.8573
